$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Kontent (column D) score updates ---
$ws.Range("D10").Value = 8
$ws.Range("D11").Value = 8
$ws.Range("D16").Value = 7
$ws.Range("D22").Value = 9
$ws.Range("D28").Value = 9

# --- Sitecore (column B) score update ---
$ws.Range("B42").Value = 9

# --- View state: scroll/selection to match author's last position ---
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D29").Select()
